# Fruta / hortaliza, semanal
# Applies a cyclic re-shuffle of the per-row sample data (columns D, L, M, N,
# O, P, Q, R, S, T) across rows 2-13 of the active sheet. Row 4 is unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for each row, keyed by row number.
# Order of fields per row: D, L, M, N, O, P, Q, R, S, T
$rows = @{
    2  = @(44601, "Primera",      30,  28000, 28000, 28000, "`$/caja 18 kilos granel",   "Región de O'Higgins", 1556, 18)
    3  = @(44411, "Primera",      210, 8000,  8000,  8000,  "`$/bandeja 8 kilos",        "Región de O'Higgins", 1000, 8)
    5  = @(44264, "Calibre 100",  50,  20000, 20000, 20000, "`$/caja 18 kilos embalada", "Región de O'Higgins", 1111, 18)
    6  = @(44217, "Primera",      55,  18000, 18000, 18000, "`$/caja 18 kilos granel",   "Región de O'Higgins", 1000, 18)
    7  = @(44418, "Especial",     100, 8000,  8000,  8000,  "`$/caja 15 kilos granel",   "Región de O'Higgins", 533,  15)
    8  = @(44392, "Especial",     500, 7000,  7000,  7000,  "`$/bandeja 8 kilos",        "Región de O'Higgins", 875,  8)
    9  = @(44966, "Primera",      4,   250000,250000,250000,"`$/bins (400 kilos)",       "Región de O'Higgins", 625,  400)
    10 = @(44966, "Primera",      80,  15000, 15000, 15000, "`$/caja 18 kilos granel",   "Región de O'Higgins", 833,  18)
    11 = @(44495, "Primera",      50,  24000, 24000, 24000, "`$/bandeja 10 kilos",       "China",               2400, 10)
    12 = @(44208, "Especial",     70,  24000, 24000, 24000, "`$/caja 15 kilos granel",   "Región de O'Higgins", 1600, 15)
    # Note: P13 is intentionally left at 7000 (not 22000) to match the
    # source diff, which shows that cell's hunk as a no-op (-7000/+7000).
    13 = @(44511, "Primera",      15,  22000, 22000, 7000,  "`$/caja 15 kilos granel",   "Región de O'Higgins", 1467, 15)
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]

    $ws.Cells.Item($r, 4).Value  = $vals[0]   # D - Fecha
    $ws.Cells.Item($r, 12).Value = $vals[1]   # L - Calidad
    $ws.Cells.Item($r, 13).Value = $vals[2]   # M - Volumen
    $ws.Cells.Item($r, 14).Value = $vals[3]   # N - Precio minimo
    $ws.Cells.Item($r, 15).Value = $vals[4]   # O - Precio maximo
    $ws.Cells.Item($r, 16).Value = $vals[5]   # P - Precio promedio ponderado
    $ws.Cells.Item($r, 17).Value = $vals[6]   # Q - Unidad de comercializacion
    $ws.Cells.Item($r, 18).Value = $vals[7]   # R - Origen
    $ws.Cells.Item($r, 19).Value = $vals[8]   # S - Precio $/Kg
    $ws.Cells.Item($r, 20).Value = $vals[9]   # T - Kg / unidad
}
